$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risk Evaluation")

$ws.Range("E1").Value = "APS VW MEB"

$ws.Range("F4").Value = "178.576-15"
$ws.Range("H4").Value = "191.674-01"
$ws.Range("J4").Value = "229.847-00"
$ws.Range("L4").Value = "234.536-00"
$ws.Range("N4").Value = "323.140-00"

$ws.Range("F5").Value = "FEDER MONT-GR"
$ws.Range("H5").Value = "STECKER GEHAEUSE MONT-GR"
$ws.Range("J5").Value = "BEDIENELEMENT"
$ws.Range("L5").Value = "GEHAEUSE MONT-GR"
$ws.Range("N5").Value = "GEHAEUSE MONT-GR"

$ws.Range("F6").Value = "MET_EXTRU"
$ws.Range("J6").Value = "TECH_PLAS"
$ws.Range("L6").Value = "TECH_PLAS"
$ws.Range("N6").Value = "TECH_PLAS"

$ws.Range("H7").Value = 434166
$ws.Range("J7").Value = 412000
$ws.Range("L7").Value = 500696
$ws.Range("N7").Value = 500696

$ws.Range("H8").Value = 46.19164619164619
$ws.Range("J8").Value = 62.65356265356265
$ws.Range("L8").Value = 36.85503685503685
$ws.Range("N8").Value = 36.85503685503685

$ws.Range("F9").Value = 9
$ws.Range("H9").Value = 9
$ws.Range("J9").Value = 9
$ws.Range("L9").Value = 9
$ws.Range("N9").Value = 9

$ws.Range("J10").Value = 1111
$ws.Range("L10").Value = 1111
$ws.Range("N10").Value = 1111

$ws.Range("H11").Value = "08.2019"
$ws.Range("J11").Value = "in.Unde"
$ws.Range("L11").Value = "in.Unde"
$ws.Range("N11").Value = "08.2019"

$ws.Range("F12").Value = "10.2020"
$ws.Range("H12").Value = "10.2020"
$ws.Range("J12").Value = "10.2020"
$ws.Range("L12").Value = "10.2020"
$ws.Range("N12").Value = "10.2020"

$ws.Range("F13").Value = "BEKAERT"
$ws.Range("H13").Value = "PA66 GF40"
$ws.Range("J13").Value = "PA66 GF40"
$ws.Range("L13").Value = "PA66 GF40"
$ws.Range("N13").Value = "PA66 GF40"

$ws.Range("J14").Value = "PlzUpdate"
$ws.Range("L14").Value = "PlzUpdate"
$ws.Range("N14").Value = "PlzUpdate"

$ws.Range("F15").Value = "PlzUpdate"
$ws.Range("J15").Value = "Jerry Fang"
$ws.Range("L15").Value = "Jerry Fang"
$ws.Range("N15").Value = "Jerry Fang"

$ws.Range("J16").Value = "CHENMI20"
$ws.Range("L16").Value = "CHENMI20"
$ws.Range("N16").Value = "CHENMI20"
